# 2017-5-31 Plate 3 BCA analyses - "Add files via upload" re-edit
#
# The underlying change re-derives the blank-corrected absorbance for well
# C13 (row 42) from the table's own "Average absorbance" column instead of
# row 39, and fills in the two previously-hardcoded averages for the last
# two sample groups (rows 78 and 81) with real AVERAGE() formulas over
# their replicate wells. Downstream columns (G/I/J/K) are plain formula
# chains off of F/G, so recalculating after the edits reproduces the rest
# of the diff's cached values automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plate3")

# --- Row 42: blank-corrected absorbance now pulled from this row's own
#     "Average absorbance at 562 nm" table column instead of E39 ---
$ws.Range("F42").Formula = "=Table1[[#This Row],[Average absorbance at 562 nm]]-E33"

# --- Row 78 & 81: average absorbance across the replicate wells instead
#     of a stale copy-pasted literal ---
$ws.Range("E78").Formula = "=AVERAGE(D78:D80)"
$ws.Range("E81").Formula = "=AVERAGE(D81:D83)"

# --- View state: scrolled/zoomed out a bit further and the active cell
#     moved from J43 to J36 ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("J36").Select() | Out-Null

$excel.Calculate() | Out-Null
